# Update EPEX Spot prices workbook with the latest daily data (03-jul / 2025-07-01).
$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column T for 03-jul ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Copy the header formatting from the previous day's column (S1) onto the new
# header cell (T1) so it keeps the same bold/border/centered style, then set
# its text.
$ws1.Range("S1").Copy()
$ws1.Range("T1").PasteSpecial(-4122)
$ws1.Range("T1").Value = "03-jul"

$prixSpot = @(85, 81.62, 79, 78.93000000000001, 80.88, 86.45, 94.83, 101.77, 104.9, 90.45, 75.06, 57.35, 42.16, 34.82, 36.84, 56.67, 64.70999999999999, 81, 96.34, 109, 118.8, 105.87, 106.8, 86.56999999999999)

for ($i = 0; $i -lt $prixSpot.Length; $i++) {
    $row = 2 + $i
    $ws1.Cells.Item($row, 20).Value = $prixSpot[$i]
}

# --- Sheet "Gaz": add row 17 for 2025-07-01 ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A17").NumberFormat = "@"
$ws2.Range("A17").Value = "2025-07-01"
$ws2.Range("A17").Style = "Normal"
$ws2.Range("B17").Value = 32.95

# --- Sheet "CO2": add row 17 for 2025-07-01 ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A17").NumberFormat = "@"
$ws3.Range("A17").Value = "2025-07-01"
$ws3.Range("A17").Style = "Normal"
$ws3.Range("B17").Value = 69.36
